$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Ning%Tang%NULL%0,                 Dengju%Li%NULL%0,                 Xiong%Wang%NULL%0,                 Ziyong%Sun%NULL%0]'
$ws.Range("E3").Value = '[Yueying%Pan%NULL%1,                 Hanxiong%Guan%NULL%1,                 Shuchang%Zhou%NULL%1,                 Yujin%Wang%NULL%1,                 Qian%Li%NULL%1,                 Tingting%Zhu%NULL%1,                 Qiongjie%Hu%qjhu@outlook.com%1,                 Liming%Xia%lmxia@tjh.tjmu.edu.cn%2]'
$ws.Range("E5").Value = '[Kaiyuan%Sun%NULL%1,                 Jenny%Chen%NULL%1,                 Cécile%Viboud%NULL%1]'
$ws.Range("E6").Value = '[Yingxia%Liu%yingxialiu@hotmail.com%4,                 Yang%Yang%NULL%4,                 Cong%Zhang%NULL%0,                 Fengming%Huang%NULL%1,                 Fuxiang%Wang%NULL%1,                 Jing%Yuan%NULL%1,                 Zhaoqin%Wang%NULL%0,                 Jinxiu%Li%NULL%1,                 Jianming%Li%NULL%1,                 Cheng%Feng%NULL%1,                 Zheng%Zhang%NULL%1,                 Lifei%Wang%NULL%1,                 Ling%Peng%NULL%1,                 Li%Chen%NULL%0,                 Yuhao%Qin%NULL%1,                 Dandan%Zhao%NULL%1,                 Shuguang%Tan%NULL%1,                 Lu%Yin%NULL%2,                 Jun%Xu%NULL%1,                 Congzhao%Zhou%NULL%1,                 Chengyu%Jiang%jiang@pumc.edu.cn%1,                 Lei%Liu%liulei3322@aliyun.com%0]'
$ws.Range("E7").Value = '[Adam%Bernheim%adam.bernheim@mountsinai.org%1,                 Xueyan%Mei%NULL%4,                 Xueyan%Mei%NULL%0,                 Mingqian%Huang%NULL%4,                 Mingqian%Huang%NULL%0,                 Yang%Yang%NULL%0,                 Zahi A.%Fayad%NULL%0,                 Zahi A.%Fayad%NULL%0,                 Ning%Zhang%NULL%5,                 Ning%Zhang%NULL%0,                 Kaiyue%Diao%NULL%2,                 Kaiyue%Diao%NULL%0,                 Bin%Lin%NULL%1,                 Xiqi%Zhu%NULL%2,                 Xiqi%Zhu%NULL%0,                 Kunwei%Li%NULL%2,                 Shaolin%Li%NULL%2,                 Hong%Shan%NULL%0,                 Adam%Jacobi%NULL%3,                 Adam%Jacobi%NULL%0,                 Michael%Chung%NULL%1]'
$ws.Range("E10").Value = '[Yan%Bai%xref no email%1,           Lingsheng%Yao%xref no email%1,           Tao%Wei%xref no email%1,           Fei%Tian%xref no email%1,           Dong-Yan%Jin%xref no email%1,           Lijuan%Chen%xref no email%1,           Meiyun%Wang%xref no email%1]'
$ws.Range("E11").Value = '[Anup%Bastola%NULL%1,                 Ranjit%Sah%ranjitsah@iom.edu.np%0,                 Alfonso J%Rodriguez-Morales%NULL%1,                 Bibek Kumar%Lal%NULL%1,                 Runa%Jha%NULL%1,                 Hemant Chanda%Ojha%NULL%1,                 Bikesh%Shrestha%NULL%1,                 Daniel K W%Chu%NULL%1,                 Leo L M%Poon%NULL%1,                 Anthony%Costello%NULL%1,                 Kouichi%Morita%NULL%1,                 Basu Dev%Pandey%NULL%1]'
$ws.Range("E13").Value = '[De%Chang%xref no email%0,           Minggui%Lin%xref no email%0,           Lai%Wei%xref no email%0,           Lixin%Xie%xref no email%0,           Guangfa%Zhu%xref no email%0,           Charles S.%Dela Cruz%xref no email%0,           Lokesh%Sharma%xref no email%0]'
$ws.Range("E15").Value = '[Nanshan%Chen%NULL%0,                 Min%Zhou%NULL%0,                 Xuan%Dong%NULL%0,                 Jieming%Qu%NULL%0,                 Fengyun%Gong%NULL%0,                 Yang%Han%NULL%0,                 Yang%Qiu%NULL%0,                 Jingli%Wang%NULL%0,                 Ying%Liu%NULL%0,                 Yuan%Wei%NULL%0,                 Jia''an%Xia%NULL%0,                 Ting%Yu%NULL%0,                 Xinxin%Zhang%NULL%0,                 Li%Zhang%NULL%0]'
$ws.Range("E16").Value = '[Michael%Chung%NULL%1,                 Adam%Bernheim%NULL%2,                 Adam%Bernheim%NULL%0,                 Xueyan%Mei%NULL%0,                 Xueyan%Mei%NULL%0,                 Ning%Zhang%NULL%0,                 Ning%Zhang%NULL%0,                 Mingqian%Huang%NULL%0,                 Mingqian%Huang%NULL%0,                 Xianjun%Zeng%NULL%1,                 Jiufa%Cui%NULL%1,                 Wenjian%Xu%NULL%2,                 Wenjian%Xu%NULL%0,                 Yang%Yang%NULL%0,                 Zahi A.%Fayad%NULL%0,                 Zahi A.%Fayad%NULL%0,                 Adam%Jacobi%NULL%0,                 Kunwei%Li%NULL%0,                 Shaolin%Li%NULL%0,                 Hong%Shan%Shanhong@mail.sysu.edu.cn%0]'
$ws.Range("E17").Value = '[Ya-ni%Duan%NULL%1,                 Jie%Qin%qinjie@mail.sysu.edu.cn%2,                 Jie%Qin%qinjie@mail.sysu.edu.cn%0]'
$ws.Range("E18").Value = '[X%Fang%NULL%1,                 M%Zhao%NULL%1,                 S%Li%NULL%1,                 L%Yang%NULL%1,                 B%Wu%bingwu69@163.com%1]'
$ws.Range("E19").Value = '[Yicheng%Fang%NULL%0,                 Huangqi%Zhang%NULL%4,                 Huangqi%Zhang%NULL%0,                 Jicheng%Xie%NULL%4,                 Jicheng%Xie%NULL%0,                 Minjie%Lin%NULL%2,                 Minjie%Lin%NULL%0,                 Lingjun%Ying%NULL%2,                 Lingjun%Ying%NULL%0,                 Peipei%Pang%NULL%4,                 Peipei%Pang%NULL%0,                 Wenbin%Ji%1224190004@qq.com%4,                 Wenbin%Ji%1224190004@qq.com%0]'
$ws.Range("E20").Value = '[Yicheng%Fang%NULL%0,                 Huangqi%Zhang%NULL%0,                 Huangqi%Zhang%NULL%0,                 Yunyu%Xu%NULL%2,                 Yunyu%Xu%NULL%0,                 Jicheng%Xie%NULL%0,                 Jicheng%Xie%NULL%0,                 Peipei%Pang%NULL%0,                 Peipei%Pang%NULL%0,                 Wenbin%Ji%1224190004@qq.com%0,                 Wenbin%Ji%1224190004@qq.com%0]'
$ws.Range("E21").Value = '[Marta%Giovanetti%NULL%1,                 Domenico%Benvenuto%NULL%1,                 Silvia%Angeletti%s.angeletti@unicampus.it%2,                 Silvia%Angeletti%s.angeletti@unicampus.it%0,                 Massimo%Ciccozzi%NULL%2,                 Massimo%Ciccozzi%NULL%0]'
$ws.Range("E22").Value = '[W.%Hao%hwdokgood@hotmail.com%1,                 M.%Li%NULL%1,                 X.%Huang%NULL%1]'
$ws.Range("E23").Value = '[Michelle L.%Holshue%NULL%1,                 Chas%DeBolt%NULL%1,                 Scott%Lindquist%NULL%1,                 Kathy H.%Lofy%NULL%1,                 John%Wiesman%NULL%1,                 Hollianne%Bruce%NULL%0,                 Christopher%Spitters%NULL%0,                 Keith%Ericson%NULL%1,                 Sara%Wilkerson%NULL%1,                 Ahmet%Tural%NULL%1,                 George%Diaz%NULL%1,                 Amanda%Cohn%NULL%2,                 LeAnne%Fox%NULL%1,                 Anita%Patel%NULL%1,                 Susan I.%Gerber%NULL%1,                 Lindsay%Kim%NULL%1,                 Suxiang%Tong%NULL%1,                 Xiaoyan%Lu%NULL%1,                 Steve%Lindstrom%NULL%1,                 Mark A.%Pallansch%NULL%1,                 William C.%Weldon%NULL%1,                 Holly M.%Biggs%NULL%1,                 Timothy M.%Uyeki%NULL%0,                 Satish K.%Pillai%NULL%1]'
$ws.Range("E24").Value = '[Chaolin%Huang%NULL%0,                 Yeming%Wang%NULL%0,                 Xingwang%Li%NULL%0,                 Lili%Ren%NULL%0,                 Jianping%Zhao%NULL%0,                 Yi%Hu%NULL%0,                 Li%Zhang%NULL%0,                 Guohui%Fan%NULL%0,                 Jiuyang%Xu%NULL%0,                 Xiaoying%Gu%NULL%0,                 Zhenshun%Cheng%NULL%0,                 Ting%Yu%NULL%0,                 Jiaan%Xia%NULL%0,                 Yuan%Wei%NULL%0,                 Wenjuan%Wu%NULL%0,                 Xuelei%Xie%NULL%0,                 Wen%Yin%NULL%0,                 Hui%Li%NULL%0,                 Min%Liu%NULL%0,                 Yan%Xiao%NULL%0,                 Hong%Gao%NULL%0,                 Li%Guo%NULL%0,                 Jungang%Xie%NULL%0,                 Guangfa%Wang%NULL%0,                 Rongmeng%Jiang%NULL%0,                 Zhancheng%Gao%NULL%0,                 Qi%Jin%NULL%0,                 Jianwei%Wang%wangjw28@163.com%0,                 Bin%Cao%caobin_ben@163.com%0]'
$ws.Range("E25").Value = '[Peikai%Huang%NULL%1,                 Tianzhu%Liu%NULL%2,                 Tianzhu%Liu%NULL%0,                 Lesheng%Huang%NULL%2,                 Lesheng%Huang%NULL%0,                 Hailong%Liu%NULL%2,                 Hailong%Liu%NULL%0,                 Ming%Lei%NULL%1,                 Wangdong%Xu%NULL%2,                 Wangdong%Xu%NULL%0,                 Xiaolu%Hu%NULL%1,                 Jun%Chen%NULL%3,                 Bo%Liu%liubogzcm@163.com%3,                 Bo%Liu%liubogzcm@163.com%0]'
$ws.Range("E26").Value = '[Huijun%Chen%NULL%0,                 Juanjuan%Guo%NULL%0,                 Chen%Wang%NULL%0,                 Fan%Luo%NULL%0,                 Xuechen%Yu%NULL%0,                 Wei%Zhang%NULL%0,                 Jiafu%Li%NULL%0,                 Dongchi%Zhao%NULL%0,                 Dan%Xu%NULL%0,                 Qing%Gong%NULL%0,                 Jing%Liao%NULL%0,                 Huixia%Yang%yanghuixia@bjmu.edu.cn%0,                 Wei%Hou%houwei@whu.edu.cn%0,                 Yuanzhen%Zhang%zhangyuanzhen@whu.edu.cn%0]'
$ws.Range("E27").Value = '[Moran%Ki%NULL%1,                 NULL%NULL%NULL%0,                 NULL%NULL%NULL%0]'
$ws.Range("E28").Value = '[Junqiang%Lei%NULL%1,                 Junfeng%Li%NULL%1,                 Xun%Li%NULL%0,                 Xiaolong%Qi%qixiaolong@vip.163.com%3]'
$ws.Range("E29").Value = '[Qun%Li%NULL%0,                 Xuhua%Guan%NULL%1,                 Peng%Wu%NULL%1,                 Xiaoye%Wang%NULL%1,                 Lei%Zhou%NULL%0,                 Yeqing%Tong%NULL%1,                 Ruiqi%Ren%NULL%1,                 Kathy S.M.%Leung%NULL%1,                 Eric H.Y.%Lau%NULL%1,                 Jessica Y.%Wong%NULL%1,                 Xuesen%Xing%NULL%1,                 Nijuan%Xiang%NULL%0,                 Yang%Wu%NULL%1,                 Chao%Li%NULL%1,                 Qi%Chen%NULL%1,                 Dan%Li%NULL%1,                 Tian%Liu%NULL%1,                 Jing%Zhao%NULL%1,                 Man%Liu%NULL%1,                 Wenxiao%Tu%NULL%1,                 Chuding%Chen%NULL%1,                 Lianmei%Jin%NULL%1,                 Rui%Yang%NULL%1,                 Qi%Wang%NULL%1,                 Suhua%Zhou%NULL%1,                 Rui%Wang%NULL%2,                 Hui%Liu%NULL%0,                 Yinbo%Luo%NULL%1,                 Yuan%Liu%NULL%1,                 Ge%Shao%NULL%1,                 Huan%Li%NULL%1,                 Zhongfa%Tao%NULL%1,                 Yang%Yang%NULL%0,                 Zhiqiang%Deng%NULL%1,                 Boxi%Liu%NULL%1,                 Zhitao%Ma%NULL%1,                 Yanping%Zhang%NULL%1,                 Guoqing%Shi%NULL%1,                 Tommy T.Y.%Lam%NULL%1,                 Joseph T.%Wu%NULL%1,                 George F.%Gao%NULL%3,                 Benjamin J.%Cowling%NULL%2,                 Bo%Yang%NULL%4,                 Bo%Yang%NULL%0,                 Gabriel M.%Leung%NULL%1,                 Zijian%Feng%NULL%0]'
$ws.Range("E30").Value = '[Xiaoqi%Lin%NULL%1,                 Zhenyu%Gong%NULL%2,                 Zhenyu%Gong%NULL%0,                 Zuke%Xiao%NULL%2,                 Zuke%Xiao%NULL%0,                 Jingliang%Xiong%NULL%1,                 Bing%Fan%NULL%2,                 Jiaqi%Liu%NULL%2,                 Jiaqi%Liu%NULL%0]'
$ws.Range("E32").Value = '[Kui%Liu%NULL%0,                 Yuan-Yuan%Fang%NULL%0,                 Yan%Deng%NULL%0,                 Wei%Liu%NULL%0,                 Mei-Fang%Wang%NULL%0,                 Jing-Ping%Ma%NULL%0,                 Wei%Xiao%NULL%0,                 Ying-Nan%Wang%NULL%0,                 Min-Hua%Zhong%NULL%0,                 Cheng-Hong%Li%NULL%0,                 Guang-Cai%Li%NULL%0,                 Hui-Guo%Liu%NULL%0,                 Xiu-Yuan%Hao%NULL%0,                 Xiu-Yuan%Hao%NULL%0,                 Pei-Fang%Wei%NULL%0]'
$ws.Range("E34").Value = '[Peng%Liu%NULL%1,                 Xian-zheng%Tan%xianzhengtan@163.com%2,                 Xian-zheng%Tan%xianzhengtan@163.com%0]'
$ws.Range("E35").Value = '[Feng%Pan%NULL%0,                 Tianhe%Ye%NULL%0,                 Tianhe%Ye%NULL%0,                 Peng%Sun%NULL%0,                 Peng%Sun%NULL%0,                 Shan%Gui%NULL%0,                 Shan%Gui%NULL%0,                 Bo%Liang%NULL%0,                 Bo%Liang%NULL%0,                 Lingli%Li%NULL%0,                 Lingli%Li%NULL%0,                 Dandan%Zheng%NULL%0,                 Dandan%Zheng%NULL%0,                 Jiazheng%Wang%NULL%0,                 Jiazheng%Wang%NULL%0,                 Richard L.%Hesketh%NULL%0,                 Richard L.%Hesketh%NULL%0,                 Lian%Yang%yanglian@hust.edu.cn%0,                 Lian%Yang%yanglian@hust.edu.cn%0,                 Chuansheng%Zheng%NULL%0,                 Chuansheng%Zheng%NULL%0]'
$ws.Range("E36").Value = '[Lan T.%Phan%NULL%1,                 Thuong V.%Nguyen%NULL%1,                 Quang C.%Luong%NULL%1,                 Thinh V.%Nguyen%NULL%1,                 Hieu T.%Nguyen%NULL%1,                 Hung Q.%Le%NULL%2,                 Hung Q.%Le%NULL%0,                 Thuc T.%Nguyen%NULL%1,                 Thang M.%Cao%NULL%2,                 Thang M.%Cao%NULL%0,                 Quang D.%Pham%NULL%1]'
$ws.Range("E37").Value = '[Wannarat A.%Pongpirul%NULL%1,                 Krit%Pongpirul%NULL%2,                 Krit%Pongpirul%NULL%0,                 Anuttra C.%Ratnarathon%NULL%2,                 Anuttra C.%Ratnarathon%NULL%0,                 Wisit%Prasithsirikul%NULL%1]'
$ws.Range("E38").Value = '[Li-Li%Ren%NULL%0,                 Ye-Ming%Wang%NULL%1,                 Zhi-Qiang%Wu%NULL%1,                 Zi-Chun%Xiang%NULL%1,                 Li%Guo%NULL%0,                 Teng%Xu%NULL%1,                 Yong-Zhong%Jiang%NULL%1,                 Yan%Xiong%NULL%0,                 Yong-Jun%Li%NULL%1,                 Xing-Wang%Li%NULL%1,                 Hui%Li%NULL%0,                 Guo-Hui%Fan%NULL%1,                 Xiao-Ying%Gu%NULL%1,                 Yan%Xiao%NULL%0,                 Hong%Gao%NULL%0,                 Jiu-Yang%Xu%NULL%1,                 Fan%Yang%NULL%1,                 Xin-Ming%Wang%NULL%1,                 Chao%Wu%NULL%1,                 Lan%Chen%NULL%1,                 Yi-Wei%Liu%NULL%1,                 Bo%Liu%NULL%0,                 Jian%Yang%NULL%1,                 Xiao-Rui%Wang%NULL%1,                 Jie%Dong%NULL%1,                 Li%Li%NULL%0,                 Chao-Lin%Huang%NULL%1,                 Jian-Ping%Zhao%NULL%1,                 Yi%Hu%NULL%0,                 Zhen-Shun%Cheng%NULL%1,                 Lin-Lin%Liu%NULL%1,                 Zhao-Hui%Qian%NULL%1,                 Chuan%Qin%NULL%1,                 Qi%Jin%NULL%0,                 Bin%Cao%NULL%0,                 Jian-Wei%Wang%NULL%1,                 Xiu-Yuan%Hao%NULL%0,                 Xiu-Yuan%Hao%NULL%0,                 Pei-Fang%Wei%NULL%0]'
$ws.Range("E39").Value = '[Camilla%Rothe%NULL%1,                 Mirjam%Schunk%NULL%1,                 Peter%Sothmann%NULL%1,                 Gisela%Bretzel%NULL%1,                 Guenter%Froeschl%NULL%1,                 Claudia%Wallrauch%NULL%1,                 Thorbjörn%Zimmer%NULL%1,                 Verena%Thiel%NULL%1,                 Christian%Janke%NULL%1,                 Wolfgang%Guggemos%NULL%2,                 Wolfgang%Guggemos%NULL%0,                 Michael%Seilmaier%NULL%1,                 Christian%Drosten%NULL%2,                 Christian%Drosten%NULL%0,                 Patrick%Vollmar%NULL%2,                 Patrick%Vollmar%NULL%0,                 Katrin%Zwirglmaier%NULL%1,                 Sabine%Zange%NULL%1,                 Roman%Wölfel%NULL%1,                 Michael%Hoelscher%NULL%2,                 Michael%Hoelscher%NULL%0]'
$ws.Range("E40").Value = '[Heshui%Shi%NULL%0,                 Xiaoyu%Han%NULL%2,                 Xiaoyu%Han%NULL%0,                 Chuansheng%Zheng%hqzcsxh@sina.com%0,                 Chuansheng%Zheng%hqzcsxh@sina.com%0]'
$ws.Range("E41").Value = '[William Kyle%Silverstein%NULL%1,                 Lynfa%Stroud%NULL%1,                 Graham Edward%Cleghorn%NULL%1,                 Jerome Allen%Leis%NULL%1]'
$ws.Range("E42").Value = '[Fengxiang%Song%NULL%0,                 Nannan%Shi%NULL%2,                 Nannan%Shi%NULL%0,                 Fei%Shan%NULL%1,                 Zhiyong%Zhang%NULL%1,                 Jie%Shen%NULL%1,                 Hongzhou%Lu%NULL%0,                 Yun%Ling%NULL%0,                 Yebin%Jiang%NULL%2,                 Yebin%Jiang%NULL%0,                 Yuxin%Shi%shiyuxin@shphc.org.cn%0,                 Yuxin%Shi%shiyuxin@shphc.org.cn%0]'
$ws.Range("E43").Value = '[Weier%Wang%NULL%1,                 Jianming%Tang%15900792812@163.com%1,                 Fangqiang%Wei%wdfwfq@126.com%1]'
$ws.Range("E44").Value = '[Min%Wei%xref no email%1,           Jingping%Yuan%xref no email%1,           Yu%Liu%xref no email%1,           Tao%Fu%xref no email%1,           Xue%Yu%xref no email%1,           Zhi-Jiang%Zhang%xref no email%1]'
$ws.Range("E45").Value = '[Wei-jie%Guan%NULL%0,                 Zheng-yi%Ni%NULL%0,                 Zheng-yi%Ni%NULL%0,                 Yu%Hu%NULL%0,                 Wen-hua%Liang%NULL%0,                 Chun-quan%Ou%NULL%0,                 Jian-xing%He%NULL%0,                 Lei%Liu%NULL%0,                 Hong%Shan%NULL%0,                 Chun-liang%Lei%NULL%0,                 David S.C.%Hui%NULL%0,                 Bin%Du%NULL%0,                 Lan-juan%Li%NULL%0,                 Guang%Zeng%NULL%0,                 Kwok-Yung%Yuen%NULL%0,                 Ru-chong%Chen%NULL%0,                 Chun-li%Tang%NULL%0,                 Tao%Wang%NULL%0,                 Ping-yan%Chen%NULL%0,                 Jie%Xiang%NULL%0,                 Shi-yue%Li%NULL%0,                 Jin-lin%Wang%NULL%0,                 Zi-jing%Liang%NULL%0,                 Yi-xiang%Peng%NULL%0,                 Li%Wei%NULL%0,                 Yong%Liu%NULL%0,                 Ya-hua%Hu%NULL%0,                 Peng%Peng%NULL%0,                 Jian-ming%Wang%NULL%0,                 Ji-yang%Liu%NULL%0,                 Zhong%Chen%NULL%0,                 Gang%Li%NULL%0,                 Zhi-jian%Zheng%NULL%0,                 Shao-qin%Qiu%NULL%0,                 Jie%Luo%NULL%0,                 Chang-jiang%Ye%NULL%0,                 Shao-yong%Zhu%NULL%0,                 Nan-shan%Zhong%NULL%0]'
$ws.Range("E46").Value = '[Xingzhi%Xie%NULL%2,                 Zheng%Zhong%NULL%2,                 Zheng%Zhong%NULL%0,                 Wei%Zhao%NULL%2,                 Wei%Zhao%NULL%0,                 Chao%Zheng%NULL%2,                 Chao%Zheng%NULL%0,                 Fei%Wang%NULL%2,                 Fei%Wang%NULL%0,                 Jun%Liu%junliu123@csu.edu.cn%0]'
$ws.Range("E47").Value = '[Xi%Xu%736461654@qq.com%0,                 Chengcheng%Yu%1515185140@qq.com%0,                 Lieguang%Zhang%zhlieguang@126.com%0,                 Liangping%Luo%tluolp@jnu.edu.cn%0,                 Jinxin%Liu%Liujx83710378@126.com%0]'
$ws.Range("E48").Value = '[Xiao-Wei%Xu%NULL%0,                 Xiao-Xin%Wu%NULL%0,                 Xian-Gao%Jiang%NULL%0,                 Kai-Jin%Xu%NULL%0,                 Ling-Jun%Ying%NULL%0,                 Chun-Lian%Ma%NULL%0,                 Shi-Bo%Li%NULL%0,                 Hua-Ying%Wang%NULL%0,                 Sheng%Zhang%NULL%0,                 Hai-Nv%Gao%NULL%0,                 Ji-Fang%Sheng%NULL%0,                 Hong-Liu%Cai%NULL%0,                 Yun-Qing%Qiu%NULL%0,                 Lan-Juan%Li%NULL%0]'
$ws.Range("E50").Value = '[Akdis%Cezmi A%coreGivesNoEmail%0,               Cao%Yi-yuan%coreGivesNoEmail%0,               Dong%Xiang%coreGivesNoEmail%0,               Gao%Ya-dong%coreGivesNoEmail%0,               Yan%You-qin%coreGivesNoEmail%0,               Yang%Yi-bin%coreGivesNoEmail%0,               Yuan%Ya-dong%coreGivesNoEmail%0,               Zhang%Jin-jin%coreGivesNoEmail%0]'
$ws.Range("E53").Value = '[Zhan%Zhang%doctorzhang2003@163.com%1,                 Xiaochen%Li%NULL%2,                 Xiaochen%Li%NULL%0,                 Wei%Zhang%NULL%0,                 Zheng-Li%Shi%NULL%1,                 Zhishui%Zheng%NULL%1,                 Tao%Wang%NULL%0]'
$ws.Range("E54").Value = '[Wan Beom%Park%NULL%1,                 Nak-Jung%Kwon%NULL%2,                 Nak-Jung%Kwon%NULL%0,                 Su-Jin%Choi%NULL%2,                 Su-Jin%Choi%NULL%0,                 Chang Kyung%Kang%NULL%2,                 Chang Kyung%Kang%NULL%0,                 Pyoeng Gyun%Choe%NULL%4,                 Pyoeng Gyun%Choe%NULL%0,                 Jin Yong%Kim%NULL%4,                 Jin Yong%Kim%NULL%0,                 Jiyoung%Yun%NULL%2,                 Jiyoung%Yun%NULL%0,                 Gir-Won%Lee%NULL%2,                 Gir-Won%Lee%NULL%0,                 Moon-Woo%Seong%NULL%2,                 Moon-Woo%Seong%NULL%0,                 Nam Joong%Kim%NULL%2,                 Nam Joong%Kim%NULL%0,                 Jeong-Sun%Seo%NULL%2,                 Jeong-Sun%Seo%NULL%0,                 Myoung-don%Oh%NULL%4,                 Myoung-don%Oh%NULL%0]'
$ws.Range("E55").Value = '[Jaegyun%Lim%NULL%1,                 Seunghyun%Jeon%NULL%2,                 Seunghyun%Jeon%NULL%0,                 Hyun-Young%Shin%NULL%2,                 Hyun-Young%Shin%NULL%0,                 Moon Jung%Kim%NULL%2,                 Moon Jung%Kim%NULL%0,                 Yu Min%Seong%NULL%2,                 Yu Min%Seong%NULL%0,                 Wang Jun%Lee%NULL%2,                 Wang Jun%Lee%NULL%0,                 Kang-Won%Choe%NULL%2,                 Kang-Won%Choe%NULL%0,                 Yu Min%Kang%NULL%2,                 Yu Min%Kang%NULL%0,                 Baeckseung%Lee%NULL%2,                 Baeckseung%Lee%NULL%0,                 Sang-Joon%Park%NULL%2,                 Sang-Joon%Park%NULL%0]'
$ws.Range("E56").Value = '[Jin Yong%Kim%NULL%1,                 Pyoeng Gyun%Choe%NULL%0,                 Pyoeng Gyun%Choe%NULL%0,                 Yoonju%Oh%NULL%2,                 Yoonju%Oh%NULL%0,                 Kyung Joong%Oh%NULL%2,                 Kyung Joong%Oh%NULL%0,                 Jinsil%Kim%NULL%2,                 Jinsil%Kim%NULL%0,                 So Jeong%Park%NULL%2,                 So Jeong%Park%NULL%0,                 Ji Hye%Park%NULL%2,                 Ji Hye%Park%NULL%0,                 Hye Kyoung%Na%NULL%2,                 Hye Kyoung%Na%NULL%0,                 Myoung-don%Oh%NULL%0,                 Myoung-don%Oh%NULL%0]'
$ws.Range("E57").Value = '[ Y.%Yang%null%1,                  Q.% Lu%null%1,                  M.% Liu%null%1,                  Y.% Wang%null%1,                  A.% Zhang%null%1,                  N.% Jalali%null%1,                  N.% Dean%null%1,                  I.% Longini%null%1,                  M. E.% Halloran%null%1,                  B.% Xu%null%1,                  X.% Zhang%null%1,                  L.% Wang%null%1,                  W.% Liu%null%1,                  L. % Fang%null%1]'
$ws.Range("E58").Value = '[Na%Zhu%NULL%3,                 Dingyu%Zhang%NULL%3,                 Wenling%Wang%NULL%2,                 Xingwang%Li%NULL%0,                 Bo%Yang%NULL%0,                 Jingdong%Song%NULL%2,                 Xiang%Zhao%NULL%2,                 Baoying%Huang%NULL%2,                 Weifeng%Shi%NULL%2,                 Roujian%Lu%NULL%2,                 Peihua%Niu%NULL%2,                 Faxian%Zhan%NULL%2,                 Xuejun%Ma%NULL%2,                 Dayan%Wang%NULL%2,                 Wenbo%Xu%NULL%2,                 Guizhen%Wu%NULL%2,                 George F.%Gao%NULL%0,                 Wenjie%Tan%NULL%2]'
$ws.Range("E59").Value = '[Jasper Fuk-Woo%Chan%NULL%0,                 Shuofeng%Yuan%NULL%1,                 Kin-Hang%Kok%NULL%1,                 Kelvin Kai-Wang%To%NULL%1,                 Hin%Chu%NULL%1,                 Jin%Yang%NULL%1,                 Fanfan%Xing%NULL%1,                 Jieling%Liu%NULL%1,                 Cyril Chik-Yan%Yip%NULL%1,                 Rosana Wing-Shan%Poon%NULL%1,                 Hoi-Wah%Tsoi%NULL%1,                 Simon Kam-Fai%Lo%NULL%1,                 Kwok-Hung%Chan%NULL%0,                 Vincent Kwok-Man%Poon%NULL%1,                 Wan-Mui%Chan%NULL%1,                 Jonathan Daniel%Ip%NULL%1,                 Jian-Piao%Cai%NULL%1,                 Vincent Chi-Chung%Cheng%NULL%1,                 Honglin%Chen%NULL%1,                 Christopher Kim-Ming%Hui%NULL%1,                 Kwok-Yung%Yuen%NULL%0]'
